$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the implication matrix values (rows 8-10, columns D and E)
$ws.Range("D8").Value = 0.5
$ws.Range("E8").Value = 0.5

$ws.Range("D9").Value = 0.6
$ws.Range("E9").Value = 0.4

$ws.Range("D10").Value = 0.8
$ws.Range("E10").Value = 0.2

# Update the active selection to F9
$ws.Range("F9").Select()
